# Add a "Postal code" (Почтовый индекс) column to the receipts template.
# This inserts a new column C (shifting the old C/D/E columns to D/E/F),
# fills in the new header cell + numbering, merges the two-row headers
# the same way the existing "Наименование..." header is merged, and
# restores the distinctive column width / selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C; everything from C onward
# (including merged ranges, styles and formulas) shifts one column right.
$ws.Columns("C").Insert() | Out-Null

# Populate the new header cell and make it wrap like the neighbouring
# "№ счета-извещения..." header; also mark the second header row cell
# so the eventual 2-row merge keeps the wrap formatting on both halves.
$ws.Range("C10").Value = "Почтовый индекс"
$ws.Range("C10").WrapText = $true
$ws.Range("C11").WrapText = $true

# Vertically merge the two-row header cells in columns B and C, matching
# how "Наименование населенного пункта / улицы" already spans B10:B11.
$ws.Range("B10:B11").Merge() | Out-Null
$ws.Range("C10:C11").Merge() | Out-Null

# Give the new column its own width (close to the original template).
$ws.Columns("C").ColumnWidth = 13.3

# Row 12 becomes the running column-number row; besides the first
# ("1", kept as text) the rest are plain numbers 2-5 under the new layout.
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 5

# Leave the same cell selected/active as in the authored workbook.
$ws.Range("C10:C11").Select() | Out-Null
